# Add a new "BB" column (next period) to the YoY rate-of-change series.
# - BB1 gets the new period's end-date (date serial 45986), formatted like BA1.
# - For every data row that already has a value in column BA, column BB is
#   populated with the same (carried-forward/latest) value as BA for that row.
# - Rows 2 and 22 have no BA value, so they stay without a BB value too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell BB1: copy BA1's formatting (date number format, font, border,
#     alignment) then set the new date value ---
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# --- Data rows: duplicate BA's value into BB for every row that has one ---
$dataRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)

foreach ($r in $dataRows) {
    $baCell = $ws.Cells.Item($r, 53)   # column BA = 53
    $bbCell = $ws.Cells.Item($r, 54)   # column BB = 54
    $bbCell.Value = $baCell.Value2
}
